$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format column D as Text to preserve literal numeric-looking strings
# (values like "1.00", "14.50", "5.30" must stay as text, matching the source data)
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "37.151.62"
$ws.Range("E2").Value = "  +0.73%  "

$ws.Range("D3").Value = "2.047.76"
$ws.Range("E3").Value = "  -3.64%  "

$ws.Range("D5").Value = "249.64"
$ws.Range("E5").Value = "  -2.91%  "

$ws.Range("D6").Value = "0.659"
$ws.Range("E6").Value = "  -2.11%  "

$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("D8").Value = "56.21"
$ws.Range("E8").Value = "  +19.99%  "

$ws.Range("D9").Value = "62.11"
$ws.Range("E9").Value = "  +0.54%  "

$ws.Range("D10").Value = "0.381"
$ws.Range("E10").Value = "  +1.36%  "

$ws.Range("D11").Value = "0.0763"
$ws.Range("E11").Value = "  +2.29%  "

$ws.Range("E12").Value = "  +5.52%  "

$ws.Range("D13").Value = "15.20"
$ws.Range("E13").Value = "  +4.11%  "

$ws.Range("D14").Value = "2.342.50"
$ws.Range("E14").Value = "  -3.68%  "

$ws.Range("D15").Value = "0.830"
$ws.Range("E15").Value = "  -3.26%  "

$ws.Range("D16").Value = "5.30"
$ws.Range("E16").Value = "  +2.09%  "

$ws.Range("D17").Value = "2.054.59"
$ws.Range("E17").Value = "  -3.30%  "

$ws.Range("D18").Value = "37.087.24"
$ws.Range("E18").Value = "  +0.43%  "

$ws.Range("D19").Value = "72.69"
$ws.Range("E19").Value = "  -2.59%  "

$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").Value = "0.0₃0880"
$ws.Range("E20").Value = "  +3.61%  "

$ws.Range("B21").Value = "Avalanche"
$ws.Range("C21").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D21").Value = "14.50"
$ws.Range("E21").Value = "  +7.13%  "

$ws.Range("D22").Value = "238.77"
$ws.Range("E22").Value = "  -1.73%  "

$ws.Range("D23").Value = "5.29"
$ws.Range("E23").Value = "  +0.48%  "

$ws.Range("E24").Value = "  +0.02%  "

$ws.Range("D25").Value = "2.43"
$ws.Range("E25").Value = "  -1.89%  "

$ws.Range("D26").Value = "170.81"
$ws.Range("E26").Value = "  -1.56%  "

$ws.Range("D27").Value = "9.23"
$ws.Range("E27").Value = "  -0.68%  "

$ws.Range("D28").Value = "20.44"
$ws.Range("E28").Value = "  -4.91%  "

$ws.Range("D29").Value = "2.02"
$ws.Range("E29").Value = "  -1.34%  "

$ws.Range("E30").Value = "  -0.73%  "

$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").Value = "4.60"
$ws.Range("E31").Value = "  +0.33%  "

$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").Value = "1.07"
$ws.Range("E32").Value = "  +16.66%  "

$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").Value = "0.0633"
$ws.Range("E33").Value = "  +4.83%  "

$ws.Range("B34").Value = "Gas"
$ws.Range("C34").Value = "https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas"
$ws.Range("D34").Value = "21.90"
$ws.Range("E34").Value = "  -5.53%  "

$ws.Range("D35").Value = "4.41"
$ws.Range("E35").Value = "  +4.10%  "

$ws.Range("E36").Value = "  +0.06%  "

$ws.Range("E37").Value = "  -5.13%  "

$ws.Range("D38").Value = "0.0859"
$ws.Range("E38").Value = "  -11.40%  "

$ws.Range("D39").Value = "1.78"
$ws.Range("E39").Value = "  -5.24%  "

$ws.Range("B40").Value = "Cronos"
$ws.Range("C40").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D40").Value = "0.109"
$ws.Range("E40").Value = "  +29.34%  "

$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").Value = "1.36"
$ws.Range("E41").Value = "  -1.16%  "

$ws.Range("D42").Value = "18.34"
$ws.Range("E42").Value = "  +11.98%  "

$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").Value = "0.0226"
$ws.Range("E43").Value = "  +0.61%  "

$ws.Range("B44").Value = "ARBITRUM"
$ws.Range("C44").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D44").Value = "1.15"
$ws.Range("E44").Value = "  -4.73%  "

$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").Value = "97.96"
$ws.Range("E45").Value = "  -2.02%  "

$ws.Range("B46").Value = "FTXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D46").Value = "4.37"
$ws.Range("E46").Value = "  +65.88%  "

$ws.Range("D47").Value = "2.79"
$ws.Range("E47").Value = "  -1.68%  "

$ws.Range("D48").Value = "1.309.01"
$ws.Range("E48").Value = "  -4.48%  "

$ws.Range("E49").Value = "  +2.91%  "

$ws.Range("D50").Value = "2.91"
$ws.Range("E50").Value = "  +2.27%  "

$ws.Range("D51").Value = "6.90"
$ws.Range("E51").Value = "  -1.16%  "
